$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Calibration bug fix: correct the protein sequence typo (STES... -> SCES...) ---
# This shared string is used by E2, E3 and E4; updating all three retires the old
# string and the corrected text becomes the new shared-string entry.
$correctedSeq = "SCESALSYAALILADSEIEISSEKLLTLTNAANVPVENIWADIFAKALDGQNLKDLLVNFSAGAAAPAGVAGGVAGGEAGEAEAEKEEEEAKEESDDDMGFGLFD"
$ws.Range("E2").Value = $correctedSeq
$ws.Range("E3").Value = $correctedSeq
$ws.Range("E4").Value = $correctedSeq

# --- Calibration bug fix: corrected observed precursor mass (neucode td) ---
$ws.Range("Q2").Value = 10894.13
$ws.Range("Q3").Value = 10894.13
$ws.Range("Q4").Value = 10894.13

# --- Widen column E so the long sequence text is readable, and set column J width ---
$ws.Columns("E").ColumnWidth = 193.83333333333334
$ws.Columns("J").ColumnWidth = 8.333333333333334

# --- Scroll the view so column F is left-most and select Q2 ---
$excel.ActiveWindow.ScrollColumn = 6
$ws.Range("Q2").Select()
